$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Third Iteration")
$ws.Activate()

# New row 15 values (set in column order so new shared strings are
# interned in the same order as the target workbook: A, E, H, I introduce
# new strings; B/C/D/F/G/M reuse "*"; J/N/P reuse existing strings)
$ws.Range("A15").Value = "Dead_sensor_alert"
$ws.Range("B15").Value = "*"
$ws.Range("C15").Value = "*"
$ws.Range("D15").Value = "*"
$ws.Range("E15").Value = "<alias> has not reported recently."
$ws.Range("F15").Value = "*"
$ws.Range("G15").Value = "*"
$ws.Range("H15").Value = "*_LATEST"
$ws.Range("I15").Value = "<now> - UTCDateTime"
$ws.Range("J15").Value = "UTCDateTime"
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = "*"
$ws.Range("N15").Value = ">"
$ws.Range("O15").Value = 1
$ws.Range("P15").Value = "alert"

# New header cell (interned right after the row-15 strings, before the
# final "subtraction..." comment string)
$ws.Range("Q1").Value = "comments"

$ws.Range("Q15").Value = "subtraction and alert value are in hours"

# Match the author's scrolled/selected view (scrolled so column P is at the
# left edge, with AD7 as the active cell)
$excel.ActiveWindow.ScrollColumn = 16
$ws.Range("AD7").Select()
